$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5 (all columns A:AH)
$ws.Range("A2").Value2 = 45145.50694444445
$ws.Range("B2").Value2 = 17.295
$ws.Range("C2").Value2 = 11.407
$ws.Range("D2").Value2 = 3.988
$ws.Range("E2").Value2 = 36.822
$ws.Range("F2").Value2 = 28.934
$ws.Range("G2").Value2 = 13.611
$ws.Range("H2").Value2 = 42.104
$ws.Range("I2").Value2 = 20.942
$ws.Range("J2").Value2 = 8.632
$ws.Range("K2").Value2 = 12.848
$ws.Range("L2").Value2 = 14.482
$ws.Range("M2").Value2 = 15.027
$ws.Range("N2").Value2 = 4.344
$ws.Range("O2").Value2 = 13.535
$ws.Range("P2").Value2 = 18.799
$ws.Range("Q2").Value2 = 11.717
$ws.Range("R2").Value2 = 3.386
$ws.Range("S2").Value2 = 2.226
$ws.Range("T2").Value2 = 198.649
$ws.Range("U2").Value2 = 37.584
$ws.Range("V2").Value2 = 12.493
$ws.Range("W2").Value2 = 24.537
$ws.Range("X2").Value2 = 12.435
$ws.Range("Y2").Value2 = 3.148
$ws.Range("Z2").Value2 = 21.511
$ws.Range("AA2").Value2 = 11.035
$ws.Range("AB2").Value2 = 10.064
$ws.Range("AC2").Value2 = 11.833
$ws.Range("AD2").Value2 = 15.025
$ws.Range("AE2").Value2 = 3.317
$ws.Range("AF2").Value2 = 37.614
$ws.Range("AG2").Value2 = 6.758
$ws.Range("AH2").Value2 = 15.619

$ws.Range("A3").Value2 = 45145.51388888889
$ws.Range("B3").Value2 = 24.021
$ws.Range("C3").Value2 = 17.322
$ws.Range("D3").Value2 = 2.066
$ws.Range("E3").Value2 = 52.091
$ws.Range("F3").Value2 = 42.151
$ws.Range("G3").Value2 = 18.904
$ws.Range("H3").Value2 = 71.25700000000001
$ws.Range("I3").Value2 = 29.086
$ws.Range("J3").Value2 = 12.725
$ws.Range("K3").Value2 = 18.77
$ws.Range("L3").Value2 = 20.857
$ws.Range("M3").Value2 = 21.894
$ws.Range("N3").Value2 = 6.038
$ws.Range("O3").Value2 = 18.798
$ws.Range("P3").Value2 = 26.609
$ws.Range("Q3").Value2 = 16.02
$ws.Range("R3").Value2 = 1.594
$ws.Range("S3").Value2 = 1.31
$ws.Range("T3").Value2 = 278.818
$ws.Range("U3").Value2 = 52.56
$ws.Range("V3").Value2 = 17.351
$ws.Range("W3").Value2 = 35.068
$ws.Range("X3").Value2 = 18.275
$ws.Range("Y3").Value2 = 3.131
$ws.Range("Z3").Value2 = 35.042
$ws.Range("AA3").Value2 = 15.326
$ws.Range("AB3").Value2 = 13.724
$ws.Range("AC3").Value2 = 16.121
$ws.Range("AD3").Value2 = 21.804
$ws.Range("AE3").Value2 = 1.247
$ws.Range("AF3").Value2 = 64.971
$ws.Range("AG3").Value2 = 9.66
$ws.Range("AH3").Value2 = 21.693

$ws.Range("A4").Value2 = 45145.52083333334
$ws.Range("B4").Value2 = 24.502
$ws.Range("C4").Value2 = 17.939
$ws.Range("D4").Value2 = 1.593
$ws.Range("E4").Value2 = 53.231
$ws.Range("F4").Value2 = 43.339
$ws.Range("G4").Value2 = 19.282
$ws.Range("H4").Value2 = 75.22199999999999
$ws.Range("I4").Value2 = 29.668
$ws.Range("J4").Value2 = 13.113
$ws.Range("K4").Value2 = 19.365
$ws.Range("L4").Value2 = 21.351
$ws.Range("M4").Value2 = 22.478
$ws.Range("N4").Value2 = 6.159
$ws.Range("O4").Value2 = 19.174
$ws.Range("P4").Value2 = 27.241
$ws.Range("Q4").Value2 = 16.227
$ws.Range("R4").Value2 = 1.105
$ws.Range("S4").Value2 = 1.069
$ws.Range("T4").Value2 = 284.545
$ws.Range("U4").Value2 = 53.639
$ws.Range("V4").Value2 = 17.698
$ws.Range("W4").Value2 = 35.959
$ws.Range("X4").Value2 = 18.808
$ws.Range("Y4").Value2 = 2.956
$ws.Range("Z4").Value2 = 36.475
$ws.Range("AA4").Value2 = 15.633
$ws.Range("AB4").Value2 = 13.912
$ws.Range("AC4").Value2 = 16.348
$ws.Range("AD4").Value2 = 22.373
$ws.Range("AE4").Value2 = 0.766
$ws.Range("AF4").Value2 = 68.42700000000001
$ws.Range("AG4").Value2 = 9.92
$ws.Range("AH4").Value2 = 22.127

$ws.Range("A5").Value2 = 45145.52777777778
$ws.Range("B5").Value2 = 22.58
$ws.Range("C5").Value2 = 16.62
$ws.Range("D5").Value2 = 1.31
$ws.Range("E5").Value2 = 49.08
$ws.Range("F5").Value2 = 40.05
$ws.Range("G5").Value2 = 17.77
$ws.Range("H5").Value2 = 70.01000000000001
$ws.Range("I5").Value2 = 27.34
$ws.Range("J5").Value2 = 12.13
$ws.Range("K5").Value2 = 17.92
$ws.Range("L5").Value2 = 19.7
$ws.Range("M5").Value2 = 20.76
$ws.Range("N5").Value2 = 5.68
$ws.Range("O5").Value2 = 17.67
$ws.Range("P5").Value2 = 25.14
$ws.Range("Q5").Value2 = 14.91
$ws.Range("R5").Value2 = 0.86
$ws.Range("S5").Value2 = 0.9
$ws.Range("T5").Value2 = 261.65
$ws.Range("U5").Value2 = 49.44
$ws.Range("V5").Value2 = 16.31
$ws.Range("W5").Value2 = 33.2
$ws.Range("X5").Value2 = 17.38
$ws.Range("Y5").Value2 = 2.65
$ws.Range("Z5").Value2 = 33.82
$ws.Range("AA5").Value2 = 14.41
$ws.Range("AB5").Value2 = 12.79
$ws.Range("AC5").Value2 = 15.03
$ws.Range("AD5").Value2 = 20.66
$ws.Range("AE5").Value2 = 0.55
$ws.Range("AF5").Value2 = 63.65
$ws.Range("AG5").Value2 = 9.16
$ws.Range("AH5").Value2 = 20.39

# Remove row 6 (dataset trimmed to 4 rows of data)
$ws.Rows.Item(6).Delete()

# Adjust column widths (stored OOXML width = ColumnWidth + 5/6)
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667   # -> width 8
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667   # -> width 8
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666   # -> width 9
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667   # -> width 8
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667   # -> width 8
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667   # -> width 8
